$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct a floating point rounding drift in the existing timestamp (C11)
$ws.Range("C11").Value = 45233.68721516203

# Append a new data row
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = 255
$ws.Range("C12").Value = 45234.50345805833
$ws.Range("C12").NumberFormat = $ws.Range("C11").NumberFormat
